# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
#
# Updates the "Home" row (row 2) stats on both the OFF and DEF sheets
# to reflect newly logged Week 17 data.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 232   # Short Att
$wsOff.Range("C2").Value = 169   # Short Comp
$wsOff.Range("D2").Value = 44    # Deep Att
$wsOff.Range("E2").Value = 18    # Deep Comp

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 246   # Short Att
$wsDef.Range("C2").Value = 157   # Short Comp
$wsDef.Range("D2").Value = 59    # Deep Att
$wsDef.Range("E2").Value = 30    # Deep Comp
$wsDef.Range("F2").Value = 6     # Short Int
